# Initial reels generation was done.
# Populate the "Base Reels 95.5 RTP" sheet with a full weighted reel strip
# (rows 1-2 already hold SCATTER / WILD; extend rows 3-70 with the
# weighted fruit symbols across all 5 reels, columns A:E).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Base Reels 95.5 RTP")
$ws.Activate()

$ranges = @(
    @{ From = 3;  To = 3;  Symbol = "WILD" },
    @{ From = 4;  To = 7;  Symbol = "FRUIT01" },
    @{ From = 8;  To = 14; Symbol = "FRUIT02" },
    @{ From = 15; To = 21; Symbol = "FRUIT03" },
    @{ From = 22; To = 28; Symbol = "FRUIT04" },
    @{ From = 29; To = 42; Symbol = "FRUIT05" },
    @{ From = 43; To = 56; Symbol = "FRUIT06" },
    @{ From = 57; To = 70; Symbol = "FRUIT07" }
)

foreach ($entry in $ranges) {
    for ($r = $entry.From; $r -le $entry.To; $r++) {
        $addr = "A" + $r + ":E" + $r
        $ws.Range($addr).Value = $entry.Symbol
    }
}

$ws.Range("A51").Select() | Out-Null
